$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("access")
$v = $ws.Range("E2").Value2
Write-Host $v
